# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 135 for
# "Vega Modelo de Temuco - Poroto verde", pushing the existing
# rows 135:141 down to 136:142 and filling the new row 135 with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 135:141 down to 136:142, leaving row 135 free for the new record.
$ws.Rows.Item(135).Insert()

$ws.Range("A135").Value = 10
$ws.Range("B135").Value = "Vega Modelo de Temuco"
$ws.Range("C135").Value = "La Araucanía"
$ws.Range("D135").Value = 44753
$ws.Range("E135").Value = 9
$ws.Range("F135").Value = 100112031
$ws.Range("G135").Value = "Poroto verde"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 40
$ws.Range("K135").Value = 40000
$ws.Range("L135").Value = 40000
$ws.Range("M135").Value = 40000
$ws.Range("N135").Value = "$/malla 25 kilos"
$ws.Range("O135").Value = "Provincia de Limarí"
$ws.Range("P135").Value = 1600
$ws.Range("Q135").Value = 25
$ws.Range("R135").Value = "Hortaliza"
